$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: fix accentuation in the "how" text (G2)
$ws.Range("G2").Value = "SEPARANDO AS TAREFAS DE FORMA AO GRUPO SER MAIS AUTÔNOMO DEPENDENDO MENOS DE ALGUMA PESSOA EM ESPECÍFICO"

# Row 7: clarify the risk description (B7) and how to handle it (G7)
$ws.Range("B7").Value = "FALTA DE COMUNICAÇÃO ENTRE O GRUPO"
$ws.Range("G7").Value = "COMUNICAR CASO HAJA DIFICULDADE EM ALGUM PROCESSO DO PROJETO, OU HAJA ALGUM IMPREVISTO"

# Update the view state to match the saved workbook (selection/zoom/scroll)
$win = $excel.ActiveWindow
$win.ScrollColumn = 3
$win.ScrollRow = 1
$win.Zoom = 80
$ws.Range("G23").Select()
